$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Changes")

# New "Changes" row describing the attachments field (row 22)
$ws.Range("A22").Value = 5
$ws.Range("B22").Value = "Level 0-0"
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = "Modify"
$ws.Range("E22").Value = "attachments"
$ws.Range("F22").Value = "MyFilename.txt"
$ws.Range("G22").Value = "Add this file as an attachment to this card. The full path of the file is needed if not specified in your env variables"

# Match the wrapped-text / text-format styling used by the rest of column G
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").WrapText = $true

$ws.Rows.Item(22).RowHeight = 30

$ws.Range("G22").Select()
